# "1. Re-org auditor rules 2. Remove auto audit from excel"
#
# The legend sheet had two status labels in row 2:
#   F2 = "自動檢核" (auto audit)   - greenish/theme-tinted fill
#   G2 = "人工註記" (manual note)  - yellow fill
#
# The auto-audit column is being removed entirely, and the manual-note
# label takes its place in column F (column G becomes empty again).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "自動檢核" (auto audit) cell completely - G2's old "人工註記"
# label effectively slides into F2's spot.
$ws.Range("G2").Clear()
$ws.Range("F2").Value = "人工註記"

# G2 used a solid yellow fill (RGB 255,255,0); give F2 that same look now
# that it carries the manual-note label.
$red = 255
$green = 255
$blue = 0
$ws.Range("F2").Interior.Color = $blue * 65536 + $green * 256 + $red

# Leave the selection where the edit happened.
[void]$ws.Range("F2").Select()
